# "linked lead tests to state" -- append two new job records (rows 53-54)
# to the JOBS sheet, mirroring the existing row layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JOBS")

# Columns B..K and M..O are stored as literal text in this sheet (even
# when the text looks like a number or a date, e.g. "0", "82916",
# "9/8/2022"). Force text typing with NumberFormat "@" before writing the
# value so Excel doesn't silently coerce it to a number/date, then put the
# cell style back to Normal so no stray formatting is left behind.
# (Note: this interpreter's functions don't get their own variable scope,
# so every local here uses a name that isn't reused by the caller.)
function Set-TextCell($targetSheet, $targetAddr, $textValue) {
    $textCellRange = $targetSheet.Range($targetAddr)
    $textCellRange.NumberFormat = "@"
    $textCellRange.Value = $textValue
    $textCellRange.Style = "Normal"
}

# Assigning Value = "" is treated as "clear the cell" rather than "store an
# empty string", so the empty-text cells (column I / notes, blank here)
# need the classic force-text-entry trick instead: a lone leading
# apostrophe collapses to a genuine empty string cell, same as the
# existing blank "notes" cells elsewhere in this sheet.
function Set-EmptyTextCell($targetSheet, $targetAddr) {
    $emptyTextCellRange = $targetSheet.Range($targetAddr)
    $emptyTextCellRange.Value = "'"
    $emptyTextCellRange.Style = "Normal"
}

$newJobRows = @(
    @{
        RowNum = 53
        A = 71305
        B = "O6CC675E200"
        C = "NA"
        D = "MULTI"
        E = "0"
        F = "GOOD"
        G = "N/A"
        H = "NO"
        I = ""
        J = "ravi"
        K = "9/8/2022"
        L = $false
        M = "N/A"
        N = "N/A"
        O = "NO"
    },
    @{
        RowNum = 54
        A = 71306
        B = "6M11 40X AWM D P"
        C = "21B 43024 M"
        D = "380 420 YY"
        E = "0"
        F = "GOOD"
        G = "82916"
        H = "NO"
        I = ""
        J = "ravi"
        K = "9/8/2022"
        L = $false
        M = "N/A"
        N = "N/A"
        O = "NO"
    }
)

foreach ($jobRecord in $newJobRows) {
    $jobRowNum = $jobRecord.RowNum

    # A: jobNumber - genuine number
    $ws.Range("A$jobRowNum").Value = $jobRecord.A

    # B..K: text columns
    Set-TextCell $ws "B$jobRowNum" $jobRecord.B
    Set-TextCell $ws "C$jobRowNum" $jobRecord.C
    Set-TextCell $ws "D$jobRowNum" $jobRecord.D
    Set-TextCell $ws "E$jobRowNum" $jobRecord.E
    Set-TextCell $ws "F$jobRowNum" $jobRecord.F
    Set-TextCell $ws "G$jobRowNum" $jobRecord.G
    Set-TextCell $ws "H$jobRowNum" $jobRecord.H
    Set-EmptyTextCell $ws "I$jobRowNum"
    Set-TextCell $ws "J$jobRowNum" $jobRecord.J
    Set-TextCell $ws "K$jobRowNum" $jobRecord.K

    # L: _isDeleted - boolean
    $ws.Range("L$jobRowNum").Value = $jobRecord.L

    # M..O: text columns
    Set-TextCell $ws "M$jobRowNum" $jobRecord.M
    Set-TextCell $ws "N$jobRowNum" $jobRecord.N
    Set-TextCell $ws "O$jobRowNum" $jobRecord.O
}
